$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Have" (column C) quantities for rows 2-13.
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 2
$ws.Range("C4").Value = 4
$ws.Range("C5").Value = 2
$ws.Range("C6").Value = 8
$ws.Range("C7").Value = 1
$ws.Range("C8").Value = 1
$ws.Range("C9").Value = 1
$ws.Range("C10").Value = 2
$ws.Range("C11").Value = 2
$ws.Range("C12").Value = 4
$ws.Range("C13").Value = 2

# The "Need" formulas (E3:E13) were missing the MAX(...,0) floor that E2 already
# has; bring the whole shared-formula group in line with it at once so the
# shared-formula grouping in the XML is preserved.
$ws.Range("E3:E13").Formula = "=MAX(B3-(C3+D3),0)"

# Reflect the last-selected cell in the sheet view.
[void]$ws.Range("D12").Select()
